$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting existing row 5..142 down to 6..143.
$ws.Rows.Item(5).Insert()

# Copy the (now-shifted) row 6 values into the newly inserted row 5,
# then update the date (column D) to the new period.
$ws.Range("A6:R6").Copy()
$ws.Range("A5:R5").PasteSpecial()

$ws.Range("D5").Value = 45245
